$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.709.31"
$ws.Range("E2").Value = "  -1.99%  "
$ws.Range("D3").Value = "3.103.87"
$ws.Range("E3").Value = "  -2.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "623.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.373"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.828"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +15.40%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("D10").Value = "3.102.14"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.613"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.184"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.31%  "
$ws.Range("E13").Value = "  -4.67%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.15%  "
$ws.Range("D15").Value = "89.508.41"
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "3.675.58"
$ws.Range("E17").Value = "  -2.53%  "
$ws.Range("D18").Value = "3.101.20"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000214"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "427.63"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "84.74"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.66%  "
$ws.Range("D28").Value = "3.264.69"
$ws.Range("E28").Value = "  -2.52%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.91%  "
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "511.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.08%  "
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("E37").Value = "  -4.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.85"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.52%  "
$ws.Range("E45").Value = "  +8.32%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.61%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0713"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +16.92%  "
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "160.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.708"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.55%  "
